# Insert a new data row at row 303 (pushing the existing rows 303-362 down to
# 304-363) and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting the row shifts rows 303:362 down to 304:363 and extends the
# worksheet dimension to A1:R363, matching the target state.
$ws.Rows.Item(303).Insert()

$ws.Range("A303").Value = 4
$ws.Range("B303").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C303").Value = "Los Lagos"
$ws.Range("D303").Value = 44995
$ws.Range("E303").Value = 10
$ws.Range("F303").Value = 100112032
$ws.Range("G303").Value = "Zapallo italiano"
$ws.Range("H303").Value = "Sin especificar"
$ws.Range("I303").Value = "Primera"
$ws.Range("J303").Value = 240
$ws.Range("K303").Value = 12000
$ws.Range("L303").Value = 13000
$ws.Range("M303").Value = 12500
$ws.Range("N303").Value = "$/caja 50 unidades"
$ws.Range("O303").Value = "Región de O'Higgins"
$ws.Range("P303").Value = 250
$ws.Range("Q303").Value = 50
$ws.Range("R303").Value = "Hortaliza"
